$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (first sheet) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(6, 6).Value = 2951
$ws1.Cells.Item(8, 6).Value = 1955

$ws1.Cells.Item(9, 6).Value = 320
$ws1.Cells.Item(9, 7).Value = "已售罄"

$ws1.Cells.Item(11, 6).Value = 797
$ws1.Cells.Item(13, 6).Value = 193
$ws1.Cells.Item(14, 6).Value = 405
$ws1.Cells.Item(15, 6).Value = 1132
$ws1.Cells.Item(17, 6).Value = 60
$ws1.Cells.Item(19, 6).Value = 7086
$ws1.Cells.Item(20, 6).Value = 265
$ws1.Cells.Item(21, 6).Value = 1765
$ws1.Cells.Item(25, 6).Value = 366
$ws1.Cells.Item(26, 6).Value = 290

$ws1.Cells.Item(27, 3).Value = "杭州·二次元拾梦漫展（取消）"
$ws1.Cells.Item(27, 7).Value = "不可售"

$ws1.Cells.Item(28, 6).Value = 1115
$ws1.Cells.Item(31, 6).Value = 116
$ws1.Cells.Item(35, 6).Value = 173
$ws1.Cells.Item(36, 6).Value = 5
$ws1.Cells.Item(41, 6).Value = 265

# --- Sheet "全部类型" (fourth sheet) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(9, 6).Value = 2951
$ws4.Cells.Item(11, 6).Value = 1955

$ws4.Cells.Item(12, 6).Value = 320
$ws4.Cells.Item(12, 7).Value = "已售罄"

$ws4.Cells.Item(14, 6).Value = 797
$ws4.Cells.Item(17, 6).Value = 193
$ws4.Cells.Item(18, 6).Value = 405
$ws4.Cells.Item(19, 6).Value = 1132
$ws4.Cells.Item(21, 6).Value = 60
$ws4.Cells.Item(23, 6).Value = 7086
$ws4.Cells.Item(24, 6).Value = 265
$ws4.Cells.Item(25, 6).Value = 1765
$ws4.Cells.Item(30, 6).Value = 366
$ws4.Cells.Item(31, 6).Value = 290

$ws4.Cells.Item(32, 3).Value = "杭州·二次元拾梦漫展（取消）"
$ws4.Cells.Item(32, 7).Value = "不可售"

$ws4.Cells.Item(33, 6).Value = 1115
$ws4.Cells.Item(35, 6).Value = 64
$ws4.Cells.Item(36, 6).Value = 116
$ws4.Cells.Item(39, 6).Value = 173
$ws4.Cells.Item(40, 6).Value = 5
$ws4.Cells.Item(45, 6).Value = 265
